$wb = $excel.ActiveWorkbook

# The "AppleTablets" sheet already has a cell styled with the custom
# "Menlo" font (style index 1 in styles.xml). Re-using that cell's font
# (via Copy) lets new cells pick up the very same font entry instead of
# Excel fabricating a new, slightly different one.
$appleTablets = $wb.Worksheets.Item("AppleTablets")
$fontDonor = $appleTablets.Range("A2")

# --- Add the two new worksheets at the end of the workbook ---------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNewest = $wb.Worksheets.Add($null, $last)
$wsNewest.Name = "NewestWatches"

$wsCompare = $wb.Worksheets.Add($null, $wsNewest)
$wsCompare.Name = "CompareHotspots"

# --- NewestWatches sheet ---------------------------------------------------
$wsNewest.Range("A2").NumberFormat = "@"
$wsNewest.Range("A2").Value = "1"

$fontDonor.Copy($wsNewest.Range("B2"))
$wsNewest.Range("B2").NumberFormat = "@"
$wsNewest.Range("B2").Value = "Newest"

$wsNewest.Range("A4").Select()

# --- CompareHotspots sheet --------------------------------------------------
$wsCompare.Range("A2").NumberFormat = "@"
$wsCompare.Range("A2").Value = "0"

$wsCompare.Range("B2").NumberFormat = "@"
$wsCompare.Range("B2").Value = "1"

$wsCompare.Range("C2").NumberFormat = "@"
$wsCompare.Range("C2").Value = "2"

$wsCompare.Range("D2").NumberFormat = "@"
$wsCompare.Range("D2").Value = "3"

$wsCompare.Range("E2").NumberFormat = "@"
$wsCompare.Range("E2").Value = "Compare Phones & Devices"

$wsCompare.Range("D3").Select()

# --- Workbook level view state ---------------------------------------------
# Applying the text cell-style to AppleTablets removed its "active" status,
# and CompareHotspots -- the last sheet -- becomes the active / selected tab,
# matching the source workbook's final state.
$wsCompare.Activate()
